$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.654043666666666
$ws.Range("H2").Value = 19.962131
$ws.Range("I2").Value = 0.3091924566209486
$ws.Range("J2").Value = 0.3091924566209486
$ws.Range("O2").Value = 0.1233263507762945
$ws.Range("P2").Value = 0.1233263507762945
$ws.Range("Q2").Value = 1.004700707273667
$ws.Range("R2").Value = 9.042306365463
$ws.Range("S2").Value = 0.03813157736261931
$ws.Range("T2").Value = 0.03813157736261931
$ws.Range("G3").Value = 6.654043666666666
$ws.Range("H3").Value = 19.962131
$ws.Range("I3").Value = 0.3091924566209486
$ws.Range("J3").Value = 0.3091924566209486
$ws.Range("O3").Value = 0.5380499444317692
$ws.Range("P3").Value = 0.5380499444317693
$ws.Range("Q3").Value = 4.383322431227445
$ws.Range("R3").Value = 39.449901881047
$ws.Range("S3").Value = 0.1663609841036236
$ws.Range("T3").Value = 0.1663609841036236
$ws.Range("G4").Value = 6.654043666666666
$ws.Range("H4").Value = 19.962131
$ws.Range("I4").Value = 0.3091924566209486
$ws.Range("J4").Value = 0.3091924566209486
$ws.Range("O4").Value = 0.0418994261307359
$ws.Range("P4").Value = 0.0418994261307359
$ws.Range("Q4").Value = 0.3413413500272222
$ws.Range("R4").Value = 3.072072150245
$ws.Range("S4").Value = 0.0129549864963702
$ws.Range("T4").Value = 0.0129549864963702
$ws.Range("G5").Value = 6.654043666666666
$ws.Range("H5").Value = 19.962131
$ws.Range("I5").Value = 0.3091924566209486
$ws.Range("J5").Value = 0.3091924566209486
$ws.Range("M5").Value = 0.3632856666666666
$ws.Range("N5").Value = 1.089857
$ws.Range("O5").Value = 0.2967242786612004
$ws.Range("P5").Value = 0.2967242786612004
$ws.Range("Q5").Value = 2.417318689474111
$ws.Range("R5").Value = 21.755868205267
$ws.Range("S5").Value = 0.09174490865833546
$ws.Range("T5").Value = 0.09174490865833546
$ws.Range("I6").Value = 0.09233579784218476
$ws.Range("J6").Value = 0.09233579784218476
$ws.Range("O6").Value = 0.1233263507762945
$ws.Range("P6").Value = 0.1233263507762945
$ws.Range("S6").Value = 0.01138743699389429
$ws.Range("T6").Value = 0.01138743699389429
$ws.Range("I7").Value = 0.09233579784218476
$ws.Range("J7").Value = 0.09233579784218476
$ws.Range("O7").Value = 0.5380499444317692
$ws.Range("P7").Value = 0.5380499444317693
$ws.Range("S7").Value = 0.04968127089805059
$ws.Range("T7").Value = 0.0496812708980506
$ws.Range("I8").Value = 0.09233579784218476
$ws.Range("J8").Value = 0.09233579784218476
$ws.Range("O8").Value = 0.0418994261307359
$ws.Range("P8").Value = 0.0418994261307359
$ws.Range("S8").Value = 0.003868816940911184
$ws.Range("T8").Value = 0.003868816940911184
$ws.Range("I9").Value = 0.09233579784218476
$ws.Range("J9").Value = 0.09233579784218476
$ws.Range("M9").Value = 0.3632856666666666
$ws.Range("N9").Value = 1.089857
$ws.Range("O9").Value = 0.2967242786612004
$ws.Range("P9").Value = 0.2967242786612004
$ws.Range("Q9").Value = 0.721896815565111
$ws.Range("R9").Value = 6.497071340085999
$ws.Range("S9").Value = 0.02739827300932869
$ws.Range("T9").Value = 0.02739827300932869
$ws.Range("G10").Value = 0.9593116666666667
$ws.Range("H10").Value = 2.877935
$ws.Range("I10").Value = 0.04457619242381535
$ws.Range("J10").Value = 0.04457619242381536
$ws.Range("O10").Value = 0.1233263507762945
$ws.Range("P10").Value = 0.1233263507762945
$ws.Range("Q10").Value = 0.1448474278616667
$ws.Range("R10").Value = 1.303626850755
$ws.Range("S10").Value = 0.005497419143131051
$ws.Range("T10").Value = 0.005497419143131052
$ws.Range("G11").Value = 0.9593116666666667
$ws.Range("H11").Value = 2.877935
$ws.Range("I11").Value = 0.04457619242381535
$ws.Range("J11").Value = 0.04457619242381536
$ws.Range("O11").Value = 0.5380499444317692
$ws.Range("P11").Value = 0.5380499444317693
$ws.Range("Q11").Value = 0.6319424033994444
$ws.Range("R11").Value = 5.687481630595
$ws.Range("S11").Value = 0.0239842178566137
$ws.Range("T11").Value = 0.02398421785661371
$ws.Range("G12").Value = 0.9593116666666667
$ws.Range("H12").Value = 2.877935
$ws.Range("I12").Value = 0.04457619242381535
$ws.Range("J12").Value = 0.04457619242381536
$ws.Range("O12").Value = 0.0418994261307359
$ws.Range("P12").Value = 0.0418994261307359
$ws.Range("Q12").Value = 0.04921108964722223
$ws.Range("R12").Value = 0.442899806825
$ws.Range("S12").Value = 0.001867716881651121
$ws.Range("T12").Value = 0.001867716881651121
$ws.Range("G13").Value = 0.9593116666666667
$ws.Range("H13").Value = 2.877935
$ws.Range("I13").Value = 0.04457619242381535
$ws.Range("J13").Value = 0.04457619242381536
$ws.Range("M13").Value = 0.3632856666666666
$ws.Range("N13").Value = 1.089857
$ws.Range("O13").Value = 0.2967242786612004
$ws.Range("P13").Value = 0.2967242786612004
$ws.Range("Q13").Value = 0.3485041783661111
$ws.Range("R13").Value = 3.136537605295
$ws.Range("S13").Value = 0.01322683854241948
$ws.Range("T13").Value = 0.01322683854241948
$ws.Range("G14").Value = 11.92023
$ws.Range("H14").Value = 35.76069
$ws.Range("I14").Value = 0.5538955531130513
$ws.Range("J14").Value = 0.5538955531130514
$ws.Range("O14").Value = 0.1233263507762945
$ws.Range("P14").Value = 0.1233263507762945
$ws.Range("Q14").Value = 1.79984744793
$ws.Range("R14").Value = 16.19862703137
$ws.Range("S14").Value = 0.0683099172766498
$ws.Range("T14").Value = 0.06830991727664981
$ws.Range("G15").Value = 11.92023
$ws.Range("H15").Value = 35.76069
$ws.Range("I15").Value = 0.5538955531130513
$ws.Range("J15").Value = 0.5538955531130514
$ws.Range("O15").Value = 0.5380499444317692
$ws.Range("P15").Value = 0.5380499444317693
$ws.Range("Q15").Value = 7.852399858169999
$ws.Range("R15").Value = 70.67159872353
$ws.Range("S15").Value = 0.2980234715734813
$ws.Range("T15").Value = 0.2980234715734815
$ws.Range("G16").Value = 11.92023
$ws.Range("H16").Value = 35.76069
$ws.Range("I16").Value = 0.5538955531130513
$ws.Range("J16").Value = 0.5538955531130514
$ws.Range("O16").Value = 0.0418994261307359
$ws.Range("P16").Value = 0.0418994261307359
$ws.Range("Q16").Value = 0.61148793195
$ws.Range("R16").Value = 5.50339138755
$ws.Range("S16").Value = 0.0232079058118034
$ws.Range("T16").Value = 0.0232079058118034
$ws.Range("G17").Value = 11.92023
$ws.Range("H17").Value = 35.76069
$ws.Range("I17").Value = 0.5538955531130513
$ws.Range("J17").Value = 0.5538955531130514
$ws.Range("M17").Value = 0.3632856666666666
$ws.Range("N17").Value = 1.089857
$ws.Range("O17").Value = 0.2967242786612004
$ws.Range("P17").Value = 0.2967242786612004
$ws.Range("Q17").Value = 4.330448702369999
$ws.Range("R17").Value = 38.97403832132999
$ws.Range("S17").Value = 0.1643542584511167
$ws.Range("T17").Value = 0.1643542584511168
